$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 29.20950566666667
$ws.Range("H2").Value = 87.628517
$ws.Range("I2").Value = 0.01829497698069002
$ws.Range("J2").Value = 0.01840828041918582
$ws.Range("M2").Value = 0.887307
$ws.Range("N2").Value = 2.661921
$ws.Range("O2").Value = 0.03991909470044044
$ws.Range("P2").Value = 0.04024932703229714
$ws.Range("Q2").Value = 25.917798844573
$ws.Range("R2").Value = 233.260189601157
$ws.Range("S2").Value = 0.0007303189186345429
$ws.Range("T2").Value = 0.000740920898694042

$ws.Range("G3").Value = 29.20950566666667
$ws.Range("H3").Value = 87.628517
$ws.Range("I3").Value = 0.01829497698069002
$ws.Range("J3").Value = 0.01840828041918582
$ws.Range("M3").Value = 0.9845453333333333
$ws.Range("O3").Value = 0.04429375447078636
$ws.Range("P3").Value = 0.04466017635322986
$ws.Range("Q3").Value = 28.75808249309022
$ws.Range("R3").Value = 258.822742437812
$ws.Range("S3").Value = 0.000810353218431372
$ws.Range("T3").Value = 0.0008221170498805469

$ws.Range("G4").Value = 29.20950566666667
$ws.Range("H4").Value = 87.628517
$ws.Range("I4").Value = 0.01829497698069002
$ws.Range("J4").Value = 0.01840828041918582
$ws.Range("M4").Value = 15.43767133333333
$ws.Range("N4").Value = 46.313014
$ws.Range("O4").Value = 0.6945260928963797
$ws.Range("P4").Value = 0.7002715882016618
$ws.Range("Q4").Value = 450.9267482911376
$ws.Range("R4").Value = 4058.340734620238
$ws.Range("S4").Value = 0.01270633888202784
$ws.Range("T4").Value = 0.01289079576520481

$ws.Range("G5").Value = 29.20950566666667
$ws.Range("H5").Value = 87.628517
$ws.Range("I5").Value = 0.01829497698069002
$ws.Range("J5").Value = 0.01840828041918582
$ws.Range("M5").Value = 0.547111
$ws.Range("N5").Value = 1.094222
$ws.Range("O5").Value = 0.02461400149063703
$ws.Range("P5").Value = 0.01654508121162658
$ws.Range("Q5").Value = 15.98084185479567
$ws.Range("R5").Value = 95.88505112877401
$ws.Range("S5").Value = 0.0004503125906738743
$ws.Range("T5").Value = 0.0003045664945018248

$ws.Range("G6").Value = 29.20950566666667
$ws.Range("H6").Value = 87.628517
$ws.Range("I6").Value = 0.01829497698069002
$ws.Range("J6").Value = 0.01840828041918582
$ws.Range("M6").Value = 4.370998666666666
$ws.Range("N6").Value = 13.112996
$ws.Range("O6").Value = 0.1966470564417564
$ws.Range("P6").Value = 0.1982738272011845
$ws.Range("Q6").Value = 127.6747103229924
$ws.Range("R6").Value = 1149.072392906932
$ws.Range("S6").Value = 0.003597653370922385
$ws.Range("T6").Value = 0.003649880210904597

$ws.Range("I7").Value = 0.913374480506715
$ws.Range("J7").Value = 0.9190311407684336
$ws.Range("M7").Value = 0.887307
$ws.Range("N7").Value = 2.661921
$ws.Range("O7").Value = 0.03991909470044044
$ws.Range("P7").Value = 0.04024932703229714
$ws.Range("Q7").Value = 1293.942926548933
$ws.Range("R7").Value = 11645.4863389404
$ws.Range("S7").Value = 0.03646108238431315
$ws.Range("T7").Value = 0.03699038493765379

$ws.Range("I8").Value = 0.913374480506715
$ws.Range("J8").Value = 0.9190311407684336
$ws.Range("M8").Value = 0.9845453333333333
$ws.Range("O8").Value = 0.04429375447078636
$ws.Range("P8").Value = 0.04466017635322986
$ws.Range("S8").Value = 0.04045678497944648
$ws.Range("T8").Value = 0.04104409282082826

$ws.Range("I9").Value = 0.913374480506715
$ws.Range("J9").Value = 0.9190311407684336
$ws.Range("M9").Value = 15.43767133333333
$ws.Range("N9").Value = 46.313014
$ws.Range("O9").Value = 0.6945260928963797
$ws.Range("P9").Value = 0.7002715882016618
$ws.Range("Q9").Value = 22512.46256837137
$ws.Range("R9").Value = 202612.1631153424
$ws.Range("S9").Value = 0.6343624092975892
$ws.Range("T9").Value = 0.643571396552696

$ws.Range("I10").Value = 0.913374480506715
$ws.Range("J10").Value = 0.9190311407684336
$ws.Range("M10").Value = 0.547111
$ws.Range("N10").Value = 1.094222
$ws.Range("O10").Value = 0.02461400149063703
$ws.Range("P10").Value = 0.01654508121162658
$ws.Range("Q10").Value = 797.8415683490756
$ws.Range("R10").Value = 4787.049410094454
$ws.Range("S10").Value = 0.0224818008247021
$ws.Range("T10").Value = 0.01520544486002755

$ws.Range("I11").Value = 0.913374480506715
$ws.Range("J11").Value = 0.9190311407684336
$ws.Range("M11").Value = 4.370998666666666
$ws.Range("N11").Value = 13.112996
$ws.Range("O11").Value = 0.1966470564417564
$ws.Range("P11").Value = 0.1982738272011845
$ws.Range("Q11").Value = 6374.144243974351
$ws.Range("R11").Value = 57367.29819576917
$ws.Range("S11").Value = 0.1796124030206639
$ws.Range("T11").Value = 0.1822198215972279

$ws.Range("G12").Value = 57.98602933333333
$ws.Range("H12").Value = 173.958088
$ws.Range("I12").Value = 0.03631876156896331
$ws.Range("J12").Value = 0.03654368891224535
$ws.Range("M12").Value = 0.887307
$ws.Range("N12").Value = 2.661921
$ws.Range("O12").Value = 0.03991909470044044
$ws.Range("P12").Value = 0.04024932703229714
$ws.Range("Q12").Value = 51.45140972967199
$ws.Range("R12").Value = 463.0626875670479
$ws.Range("S12").Value = 0.001449812082474163
$ws.Range("T12").Value = 0.001470858885995494

$ws.Range("G13").Value = 57.98602933333333
$ws.Range("H13").Value = 173.958088
$ws.Range("I13").Value = 0.03631876156896331
$ws.Range("J13").Value = 0.03654368891224535
$ws.Range("M13").Value = 0.9845453333333333
$ws.Range("O13").Value = 0.04429375447078636
$ws.Range("P13").Value = 0.04466017635322986
$ws.Range("Q13").Value = 57.0898745786631
$ws.Range("R13").Value = 513.808871207968
$ws.Range("S13").Value = 0.001608694307618693
$ws.Range("T13").Value = 0.001632047591418448

$ws.Range("G14").Value = 57.98602933333333
$ws.Range("H14").Value = 173.958088
$ws.Range("I14").Value = 0.03631876156896331
$ws.Range("J14").Value = 0.03654368891224535
$ws.Range("M14").Value = 15.43767133333333
$ws.Range("N14").Value = 46.313014
$ws.Range("O14").Value = 0.6945260928963797
$ws.Range("P14").Value = 0.7002715882016618
$ws.Range("Q14").Value = 895.1692627730256
$ws.Range("R14").Value = 8056.52336495723
$ws.Range("S14").Value = 0.02522432757132728
$ws.Range("T14").Value = 0.02559050707332551

$ws.Range("G15").Value = 57.98602933333333
$ws.Range("H15").Value = 173.958088
$ws.Range("I15").Value = 0.03631876156896331
$ws.Range("J15").Value = 0.03654368891224535
$ws.Range("M15").Value = 0.547111
$ws.Range("N15").Value = 1.094222
$ws.Range("O15").Value = 0.02461400149063703
$ws.Range("P15").Value = 0.01654508121162658
$ws.Range("Q15").Value = 31.72479449458933
$ws.Range("R15").Value = 190.348766967536
$ws.Range("S15").Value = 0.0008939500513965539
$ws.Range("T15").Value = 0.0006046183008255171

$ws.Range("G16").Value = 57.98602933333333
$ws.Range("H16").Value = 173.958088
$ws.Range("I16").Value = 0.03631876156896331
$ws.Range("J16").Value = 0.03654368891224535
$ws.Range("M16").Value = 4.370998666666666
$ws.Range("N16").Value = 13.112996
$ws.Range("O16").Value = 0.1966470564417564
$ws.Range("P16").Value = 0.1982738272011845
$ws.Range("Q16").Value = 253.4568569012942
$ws.Range("R16").Value = 2281.111712111648
$ws.Range("S16").Value = 0.007141977556146622
$ws.Range("T16").Value = 0.007245657060680376

$ws.Range("G17").Value = 29.481085
$ws.Range("H17").Value = 58.96217
$ws.Range("I17").Value = 0.01846507700595112
$ws.Range("J17").Value = 0.01238628926567028
$ws.Range("M17").Value = 0.887307
$ws.Range("N17").Value = 2.661921
$ws.Range("O17").Value = 0.03991909470044044
$ws.Range("P17").Value = 0.04024932703229714
$ws.Range("Q17").Value = 26.158773088095
$ws.Range("R17").Value = 156.95263852857
$ws.Range("S17").Value = 0.0007371091576514882
$ws.Range("T17").Value = 0.0004985398073705945

$ws.Range("G18").Value = 29.481085
$ws.Range("H18").Value = 58.96217
$ws.Range("I18").Value = 0.01846507700595112
$ws.Range("J18").Value = 0.01238628926567028
$ws.Range("M18").Value = 0.9845453333333333
$ws.Range("O18").Value = 0.04429375447078636
$ws.Range("P18").Value = 0.04466017635322986
$ws.Range("Q18").Value = 29.02546465835333
$ws.Range("R18").Value = 174.15278795012
$ws.Range("S18").Value = 0.000817887587185762
$ws.Range("T18").Value = 0.0005531738629669527

$ws.Range("G19").Value = 29.481085
$ws.Range("H19").Value = 58.96217
$ws.Range("I19").Value = 0.01846507700595112
$ws.Range("J19").Value = 0.01238628926567028
$ws.Range("M19").Value = 15.43767133333333
$ws.Range("N19").Value = 46.313014
$ws.Range("O19").Value = 0.6945260928963797
$ws.Range("P19").Value = 0.7002715882016618
$ws.Range("Q19").Value = 455.1193007800633
$ws.Range("R19").Value = 2730.71580468038
$ws.Range("S19").Value = 0.01282447778797401
$ws.Range("T19").Value = 0.00867376645599612

$ws.Range("G20").Value = 29.481085
$ws.Range("H20").Value = 58.96217
$ws.Range("I20").Value = 0.01846507700595112
$ws.Range("J20").Value = 0.01238628926567028
$ws.Range("M20").Value = 0.547111
$ws.Range("N20").Value = 1.094222
$ws.Range("O20").Value = 0.02461400149063703
$ws.Range("P20").Value = 0.01654508121162658
$ws.Range("Q20").Value = 16.129425895435
$ws.Range("R20").Value = 64.51770358174001
$ws.Range("S20").Value = 0.0004544994329492086
$ws.Range("T20").Value = 0.0002049321618112133

$ws.Range("G21").Value = 29.481085
$ws.Range("H21").Value = 58.96217
$ws.Range("I21").Value = 0.01846507700595112
$ws.Range("J21").Value = 0.01238628926567028
$ws.Range("M21").Value = 4.370998666666666
$ws.Range("N21").Value = 13.112996
$ws.Range("O21").Value = 0.1966470564417564
$ws.Range("P21").Value = 0.1982738272011845
$ws.Range("Q21").Value = 128.8617832268866
$ws.Range("R21").Value = 773.17069936132
$ws.Range("S21").Value = 0.003631103040190649
$ws.Range("T21").Value = 0.002455876977525395

$ws.Range("G22").Value = 21.628479
$ws.Range("H22").Value = 64.885437
$ws.Range("I22").Value = 0.01354670393768061
$ws.Range("J22").Value = 0.01363060063446486
$ws.Range("M22").Value = 0.887307
$ws.Range("N22").Value = 2.661921
$ws.Range("O22").Value = 0.03991909470044044
$ws.Range("P22").Value = 0.04024932703229714
$ws.Range("Q22").Value = 19.191100816053
$ws.Range("R22").Value = 172.719907344477
$ws.Range("S22").Value = 0.0005407721573671018
$ws.Range("T22").Value = 0.0005486225025832131

$ws.Range("G23").Value = 21.628479
$ws.Range("H23").Value = 64.885437
$ws.Range("I23").Value = 0.01354670393768061
$ws.Range("J23").Value = 0.01363060063446486
$ws.Range("M23").Value = 0.9845453333333333
$ws.Range("O23").Value = 0.04429375447078636
$ws.Range("P23").Value = 0.04466017635322986
$ws.Range("Q23").Value = 21.294218066548
$ws.Range("R23").Value = 191.647962598932
$ws.Range("S23").Value = 0.0006000343781040597
$ws.Range("T23").Value = 0.0006087450281356476

$ws.Range("G24").Value = 21.628479
$ws.Range("H24").Value = 64.885437
$ws.Range("I24").Value = 0.01354670393768061
$ws.Range("J24").Value = 0.01363060063446486
$ws.Range("M24").Value = 15.43767133333333
$ws.Range("N24").Value = 46.313014
$ws.Range("O24").Value = 0.6945260928963797
$ws.Range("P24").Value = 0.7002715882016618
$ws.Range("Q24").Value = 333.8933502419019
$ws.Range("R24").Value = 3005.040152177117
$ws.Range("S24").Value = 0.009408539357461316
$ws.Range("T24").Value = 0.009545122354439288

$ws.Range("G25").Value = 21.628479
$ws.Range("H25").Value = 64.885437
$ws.Range("I25").Value = 0.01354670393768061
$ws.Range("J25").Value = 0.01363060063446486
$ws.Range("M25").Value = 0.547111
$ws.Range("N25").Value = 1.094222
$ws.Range("O25").Value = 0.02461400149063703
$ws.Range("P25").Value = 0.01654508121162658
$ws.Range("Q25").Value = 11.833178774169
$ws.Range("R25").Value = 70.999072645014
$ws.Range("S25").Value = 0.0003334385909152891
$ws.Range("T25").Value = 0.0002255193944604699

$ws.Range("G26").Value = 21.628479
$ws.Range("H26").Value = 64.885437
$ws.Range("I26").Value = 0.01354670393768061
$ws.Range("J26").Value = 0.01363060063446486
$ws.Range("M26").Value = 4.370998666666666
$ws.Range("N26").Value = 13.112996
$ws.Range("O26").Value = 0.1966470564417564
$ws.Range("P26").Value = 0.1982738272011845
$ws.Range("Q26").Value = 94.53805287102799
$ws.Range("R26").Value = 850.8424758392518
$ws.Range("S26").Value = 0.002663919453832843
$ws.Range("T26").Value = 0.003649880210904597

